$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for the market date 2021-10-26 (serial 44495).
# It belongs right after the existing row for 2021-10-08 (row 10) and before the
# existing row for 2021-10-14 (row 11, which becomes row 12). Insert a new row at
# position 11 so every following row shifts down by one, preserving all of their
# existing values untouched.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new record's data.
$ws.Range("A11").Value = 12
$ws.Range("B11").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C11").Value = "Metropolitana"
$ws.Range("D11").Value = 44495
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = 100112028
$ws.Range("G11").Value = "Sandia"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 800
$ws.Range("L11").Value = 800
$ws.Range("M11").Value = 800
$ws.Range("N11").Value = '$/kilo (volumen en unidades)'
$ws.Range("O11").Value = "Perú"
$ws.Range("P11").Value = 800
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = "Hortaliza"
